$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 1.06

$ws1.Range("D3").Value = 89
$ws1.Range("H3").Value = 7.56
$ws1.Range("L3").Value = 0.87

$ws1.Range("D4").Value = 88
$ws1.Range("H4").Value = 6.6
$ws1.Range("L4").Value = 0.88

$ws1.Range("D5").Value = 89
$ws1.Range("H5").Value = 5.56
$ws1.Range("L5").Value = 0.92

$ws1.Range("D6").Value = 89
$ws1.Range("H6").Value = 4.56
$ws1.Range("L6").Value = 1.02

$ws1.Range("D7").Value = 88
$ws1.Range("H7").Value = 3.6
$ws1.Range("L7").Value = 1.18

$ws1.Range("H8").Value = 2.67
$ws1.Range("L8").Value = 1.04

$ws1.Range("D9").Value = 88
$ws1.Range("H9").Value = 1.63
$ws1.Range("L9").Value = 1.07

$ws1.Range("H10").Value = 0.64
$ws1.Range("L10").Value = 1.2

$ws1.Range("L11").Value = 0.97

$ws1.Range("L12").Value = 0.9

$ws1.Range("L13").Value = 0.97

$ws1.Range("L14").Value = 0.88

$ws1.Range("L15").Value = 0.98

$ws1.Range("L16").Value = 0.9

$ws1.Range("L17").Value = 0.8100000000000001

# --- Sheet: Summary ---
# These cells hold numeric-looking text (t="inlineStr"/shared string in the
# original file), so prefix with an apostrophe to keep them stored as text
# instead of being auto-converted to numbers by Excel.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'1317"
$ws2.Range("B10").Value = "'698"
$ws2.Range("B11").Value = "'347"
$ws2.Range("B12").Value = "'89"
